$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50049.617
$ws.Range("I6").Value = 55244.74
$ws.Range("K6").Value = 165734.22
$ws.Range("M6").Value = -165622.22
$ws.Range("H21").Value = 6000
$ws.Range("J21").Value = 6000
$ws.Range("L21").Value = 6000
$ws.Range("N21").Value = -6936
$ws.Range("H23").Value = 6000
$ws.Range("J23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("N23").Value = -6468
$ws.Range("H39").Value = 955.1
$ws.Range("I39").Value = 832.75
$ws.Range("K39").Value = 2498.25
$ws.Range("M39").Value = -2202.25
$ws.Range("H41").Value = 1130.4375
$ws.Range("I41").Value = 419.42856
$ws.Range("J41").Value = 1683.4445
$ws.Range("K41").Value = 419.42856
$ws.Range("L41").Value = 1683.4445
$ws.Range("M41").Value = 20.57144
$ws.Range("N41").Value = -2563.4445
$ws.Range("H58").Value = 629.6
$ws.Range("I58").Value = 662.25
$ws.Range("J58").Value = 499
$ws.Range("K58").Value = 1986.75
$ws.Range("L58").Value = 1497
$ws.Range("M58").Value = -1836.75
$ws.Range("N58").Value = -1797
$ws.Range("H86").Value = 79719.16
$ws.Range("I86").Value = 202559.8
$ws.Range("K86").Value = 202559.8
$ws.Range("M86").Value = -201436.8
$ws.Range("H89").Value = 79719.16
$ws.Range("I89").Value = 202559.8
$ws.Range("K89").Value = 1012799
$ws.Range("M89").Value = -1007183
$ws.Range("H97").Value = 1804.875
$ws.Range("J97").Value = 1804.875
$ws.Range("L97").Value = 5414.625
$ws.Range("N97").Value = -6406.625
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180
$ws.Range("H112").Value = 23233726
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 23233726
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 69701178
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -69703394
$ws.Range("H115").Value = 8998973
$ws.Range("I115").Value = 9997748
$ws.Range("J115").Value = 9998.5
$ws.Range("K115").Value = 29993244
$ws.Range("L115").Value = 29995.5
$ws.Range("M115").Value = -29991677
$ws.Range("N115").Value = -33129.5
$ws.Range("H132").Value = 2529.4075
$ws.Range("I132").Value = 2305.889
$ws.Range("K132").Value = 6917.667
$ws.Range("M132").Value = -4387.667
$ws.Range("H137").Value = 3065846.8
$ws.Range("I137").Value = 6742417
$ws.Range("K137").Value = 20227251
$ws.Range("M137").Value = -20224701
$ws.Range("H138").Value = 5029.91
$ws.Range("I138").Value = 2047
$ws.Range("J138").Value = 5598.0835
$ws.Range("K138").Value = 6141
$ws.Range("L138").Value = 16794.2505
$ws.Range("M138").Value = -1001
$ws.Range("N138").Value = -27074.2505
$ws.Range("H141").Value = 3763.2222
$ws.Range("I141").Value = 1514
$ws.Range("J141").Value = 6574.75
$ws.Range("K141").Value = 4542
$ws.Range("L141").Value = 19724.25
$ws.Range("M141").Value = 638
$ws.Range("N141").Value = -30084.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 95192.81
$ws.Range("I74").Value = 105086.84
$ws.Range("J74").Value = 1199.5
$ws.Range("K74").Value = 105086.84
$ws.Range("L74").Value = 1199.5
$ws.Range("M74").Value = -104212.84
$ws.Range("N74").Value = -2947.5
$ws.Range("H77").Value = 95192.81
$ws.Range("I77").Value = 105086.84
$ws.Range("J77").Value = 1199.5
$ws.Range("K77").Value = 525434.2
$ws.Range("L77").Value = 5997.5
$ws.Range("M77").Value = -521066.2
$ws.Range("N77").Value = -14733.5
$ws.Range("H92").Value = 49585.285
$ws.Range("J92").Value = 49585.285
$ws.Range("L92").Value = 49585.285
$ws.Range("N92").Value = -54577.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1913.625
$ws.Range("I64").Value = 570
$ws.Range("J64").Value = 2361.5
$ws.Range("K64").Value = 570
$ws.Range("L64").Value = 2361.5
$ws.Range("M64").Value = -345
$ws.Range("N64").Value = -2811.5
$ws.Range("H67").Value = 1913.625
$ws.Range("I67").Value = 570
$ws.Range("J67").Value = 2361.5
$ws.Range("K67").Value = 570
$ws.Range("L67").Value = 2361.5
$ws.Range("M67").Value = 210
$ws.Range("N67").Value = -3921.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 236276.4
$ws.Range("I31").Value = 372335.22
$ws.Range("K31").Value = 372335.22
$ws.Range("M31").Value = -372040.22
$ws.Range("H34").Value = 236276.4
$ws.Range("I34").Value = 372335.22
$ws.Range("K34").Value = 372335.22
$ws.Range("M34").Value = -372133.22
$ws.Range("H41").Value = 8154.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H43").Value = 22216.5
$ws.Range("J43").Value = 22216.5
$ws.Range("L43").Value = 22216.5
$ws.Range("N43").Value = -22584.5
$ws.Range("H47").Value = 4750.75
$ws.Range("I47").Value = 1031.5
$ws.Range("J47").Value = 8470
$ws.Range("K47").Value = 1031.5
$ws.Range("L47").Value = 8470
$ws.Range("M47").Value = -465.5
$ws.Range("N47").Value = -9602
$ws.Range("H62").Value = 2954.2
$ws.Range("I62").Value = 2956.5
$ws.Range("K62").Value = 2956.5
$ws.Range("M62").Value = -2332.5
$ws.Range("H65").Value = 2954.2
$ws.Range("I65").Value = 2956.5
$ws.Range("K65").Value = 14782.5
$ws.Range("M65").Value = -11662.5
$ws.Range("H86").Value = 4812.769
$ws.Range("I86").Value = 3360.2856
$ws.Range("K86").Value = 3360.2856
$ws.Range("M86").Value = -2237.2856
$ws.Range("H89").Value = 4812.769
$ws.Range("I89").Value = 3360.2856
$ws.Range("K89").Value = 16801.428
$ws.Range("M89").Value = -11185.428
$ws.Range("H101").Value = 22216.5
$ws.Range("J101").Value = 22216.5
$ws.Range("L101").Value = 22216.5
$ws.Range("N101").Value = -28706.5
$ws.Range("H134").Value = 6332.905
$ws.Range("I134").Value = 6594
$ws.Range("J134").Value = 1111
$ws.Range("K134").Value = 19782
$ws.Range("L134").Value = 3333
$ws.Range("M134").Value = -17247
$ws.Range("N134").Value = -8403

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 44873700
$ws.Range("J129").Value = 2466
$ws.Range("L129").Value = 7398
$ws.Range("N129").Value = -17398
$ws.Range("H131").Value = 5683196.5
$ws.Range("J131").Value = 1442.9012
$ws.Range("L131").Value = 4328.7036
$ws.Range("N131").Value = -14408.7036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4761.1177
$ws.Range("I80").Value = 4638.5713
$ws.Range("J80").Value = 5333
$ws.Range("K80").Value = 4638.5713
$ws.Range("L80").Value = 5333
$ws.Range("M80").Value = -3640.5713
$ws.Range("N80").Value = -7329
$ws.Range("H83").Value = 4761.1177
$ws.Range("I83").Value = 4638.5713
$ws.Range("J83").Value = 5333
$ws.Range("K83").Value = 23192.8565
$ws.Range("L83").Value = 26665
$ws.Range("M83").Value = -18200.8565
$ws.Range("N83").Value = -36649
$ws.Range("H99").Value = 29974.6
$ws.Range("I99").Value = 24958
$ws.Range("K99").Value = 24958
$ws.Range("M99").Value = -22712
$ws.Range("H102").Value = 3922.75
$ws.Range("I102").Value = 3906.5386
$ws.Range("K102").Value = 3906.5386
$ws.Range("M102").Value = -2284.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 8982
$ws.Range("J20").Value = 9959
$ws.Range("L20").Value = 9959
$ws.Range("N20").Value = -10411
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""
$ws.Range("H61").Value = 1782.091
$ws.Range("I61").Value = 1737.4445
$ws.Range("K61").Value = 1737.4445
$ws.Range("M61").Value = -1535.4445
$ws.Range("H82").Value = 1363.1666
$ws.Range("I82").Value = 731.875
$ws.Range("K82").Value = 731.875
$ws.Range("M82").Value = -370.875
$ws.Range("H85").Value = 1363.1666
$ws.Range("I85").Value = 731.875
$ws.Range("K85").Value = 731.875
$ws.Range("M85").Value = 516.125
$ws.Range("H97").Value = 23332.334
$ws.Range("J97").Value = 23332.334
$ws.Range("L97").Value = 23332.334
$ws.Range("N97").Value = -25314.334
$ws.Range("H113").Value = 1782.091
$ws.Range("I113").Value = 1737.4445
$ws.Range("K113").Value = 1737.4445
$ws.Range("M113").Value = 432.5554999999999
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("K7").Value = 4
$ws.Range("M7").Value = 109
$ws.Range("H45").Value = 36374.2
$ws.Range("J45").Value = 31312.75
$ws.Range("L45").Value = 31312.75
$ws.Range("N45").Value = -32294.75
$ws.Range("H107").Value = 1009.8333
$ws.Range("I107").Value = 1050.9
$ws.Range("J107").Value = 804.5
$ws.Range("K107").Value = 3152.7
$ws.Range("L107").Value = 2413.5
$ws.Range("M107").Value = -1232.7
$ws.Range("N107").Value = -6253.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
$ws.Range("H129").Value = 99429
$ws.Range("J129").Value = 99429
$ws.Range("L129").Value = 99429
$ws.Range("N129").Value = -109429
$ws.Range("H132").Value = 1754.9584
$ws.Range("I132").Value = 1448.1765
$ws.Range("K132").Value = 4344.529500000001
$ws.Range("M132").Value = -1814.529500000001

